$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.190.78"
$ws.Range("E2").Value = "  -0.38%  "
$ws.Range("D3").Value = "3.670.72"
$ws.Range("E3").Value = "  -0.44%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'673.81"
$ws.Range("D6").Value = "'157.23"
$ws.Range("E6").Value = "  -3.30%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -1.53%  "
$ws.Range("E9").Value = "  -1.86%  "
$ws.Range("D10").Value = "'6.95"
$ws.Range("E10").Value = "  -5.53%  "
$ws.Range("E11").Value = "  -2.44%  "
$ws.Range("E12").Value = "  -3.75%  "
$ws.Range("D13").Value = "4.289.88"
$ws.Range("E13").Value = "  -0.53%  "
$ws.Range("D14").Value = "'32.15"
$ws.Range("E14").Value = "  -4.32%  "
$ws.Range("D15").Value = "3.648.32"
$ws.Range("E15").Value = "  -1.18%  "
$ws.Range("D16").Value = "69.205.04"
$ws.Range("E16").Value = "  -0.39%  "
$ws.Range("E17").Value = "  +1.00%  "
$ws.Range("D18").Value = "'15.99"
$ws.Range("E18").Value = "  -1.51%  "
$ws.Range("D19").Value = "'6.41"
$ws.Range("E19").Value = "  -3.25%  "
$ws.Range("D20").Value = "'467.26"
$ws.Range("E20").Value = "  -3.49%  "
$ws.Range("E21").Value = "  +0.34%  "
$ws.Range("D22").Value = "'0.647"
$ws.Range("E22").Value = "  -3.12%  "
$ws.Range("D23").Value = "'79.74"
$ws.Range("E23").Value = "  -0.82%  "
$ws.Range("D24").Value = "3.816.10"
$ws.Range("E24").Value = "  -0.41%  "
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("B26").Value = "PEPE"
$ws.Range("C26").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D26").Value = "'0.0000120"
$ws.Range("E26").Value = "  -7.93%  "
$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").Value = "'10.87"
$ws.Range("E27").Value = "  -5.34%  "
$ws.Range("D28").Value = "'8.99"
$ws.Range("E28").Value = "  -6.22%  "
$ws.Range("E29").Value = "  -2.70%  "
$ws.Range("E30").Value = "  -6.24%  "
$ws.Range("D31").Value = "'6.59"
$ws.Range("E31").Value = "  -4.10%  "
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("D33").Value = "'26.79"
$ws.Range("E33").Value = "  -1.09%  "
$ws.Range("E34").Value = "  -5.74%  "
$ws.Range("D35").Value = "3.661.13"
$ws.Range("E35").Value = "  +0.08%  "
$ws.Range("D36").Value = "'0.160"
$ws.Range("E36").Value = "  -4.35%  "
$ws.Range("D37").Value = "'8.11"
$ws.Range("E37").Value = "  -4.66%  "
$ws.Range("E38").Value = "  -3.23%  "
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("E41").Value = "  -2.49%  "
$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").Value = "'173.68"
$ws.Range("E42").Value = "  +7.84%  "
$ws.Range("B43").Value = "Hedera"
$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D43").Value = "'0.0896"
$ws.Range("E43").Value = "  -4.47%  "
$ws.Range("E44").Value = "  -1.73%  "
$ws.Range("D45").Value = "'47.53"
$ws.Range("E45").Value = "  -1.92%  "
$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").Value = "'2.67"
$ws.Range("E46").Value = "  -6.59%  "
$ws.Range("B47").Value = "FLOKI"
$ws.Range("C47").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D47").Value = "'0.000275"
$ws.Range("E47").Value = "  -5.38%  "
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").Value = "'1.27"
$ws.Range("E48").Value = "  -6.27%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").Value = "'27.31"
$ws.Range("E49").Value = "  -8.75%  "
$ws.Range("E50").Value = "  -4.26%  "
$ws.Range("D51").Value = "'7.77"
$ws.Range("E51").Value = "  -3.27%  "
